# Rename the "Inputs and Outputs" sheet's title cell (A1) from "Inputs"
# to "Results Summary and Inputs".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs and Outputs")
$ws.Range("A1").Value = "Results Summary and Inputs"
$ws.Range("A1").Select() | Out-Null
